# Implementacion localDate y LocalTime en el Home
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update progress values for MenuAdministradorController (row 9) and
# MenuGerenteController (row 10) in the last week column (F) to 100.
$ws.Range("F9").Value = 100
$ws.Range("F10").Value = 100

# Move the active selection from F20 to F2.
$ws.Range("F2").Select()
